# "new changes with extended report"
#
# Updates a handful of generated test-data strings on Sheet1 (columns C/E for
# rows 2-4) and moves the active selection from I8 to E4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Order matters: new shared-string entries are appended in the order they are
# first written, so touch E2 before C2 to reproduce the same shared-strings
# ordering as the authored change.
$ws.Range("E2").Value = "16hdfc16"
$ws.Range("C2").Value = "16swati19"
$ws.Range("C3").Value = "26swati39"
$ws.Range("C4").Value = "36swati29"
$ws.Range("E4").Value = "36hdfc20"

# Move the selection/active cell from I8 to E4.
$ws.Range("E4").Select()
